$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Account Details")

# Update the password cell (B2) to hold the new password text instead of the old numeric value
$ws.Range("B2").Value = "PasswordPassword!Password!123"

# Update the active selection on the sheet
$ws.Range("D4").Select()
